$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Force column D/E cells to remain plain text (matches original inlineStr cells)
# even when the new value looks like a number (e.g. "62.062.04", "0.999").
$textRange = $ws.Range("D2:E51")
$textRange.NumberFormat = "@"

$ws.Range("D2").Value = "62.062.04"
$ws.Range("E2").Value = "  -0.27%  "

$ws.Range("D3").Value = "2.463.98"
$ws.Range("E3").Value = "  -1.68%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "550.53"
$ws.Range("E5").Value = "  -1.32%  "

$ws.Range("D6").Value = "147.07"
$ws.Range("E6").Value = "  -0.66%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "0.587"
$ws.Range("E8").Value = "  -4.05%  "

$ws.Range("D9").Value = "2.462.37"
$ws.Range("E9").Value = "  -1.71%  "

$ws.Range("E10").Value = "  -2.62%  "

$ws.Range("E11").Value = "  +0.12%  "

$ws.Range("D12").Value = "5.41"
$ws.Range("E12").Value = "  -0.95%  "

$ws.Range("D13").Value = "0.351"
$ws.Range("E13").Value = "  -3.67%  "

$ws.Range("D14").Value = "26.21"
$ws.Range("E14").Value = "  -1.39%  "

$ws.Range("D15").Value = "2.908.86"
$ws.Range("E15").Value = "  -1.77%  "

$ws.Range("E16").Value = "  +0.32%  "

$ws.Range("D17").Value = "61.963.52"
$ws.Range("E17").Value = "  -0.15%  "

$ws.Range("D18").Value = "2.463.83"
$ws.Range("E18").Value = "  -1.45%  "

$ws.Range("D19").Value = "10.96"
$ws.Range("E19").Value = "  -3.41%  "

$ws.Range("D20").Value = "7.04"
$ws.Range("E20").Value = "  -0.93%  "

$ws.Range("D21").Value = "4.17"
$ws.Range("E21").Value = "  -2.40%  "

$ws.Range("D22").Value = "321.18"
$ws.Range("E22").Value = "  -1.92%  "

$ws.Range("E23").Value = "  +0.12%  "

$ws.Range("E24").Value = "  +6.53%  "

$ws.Range("D25").Value = "64.04"
$ws.Range("E25").Value = "  -1.55%  "

$ws.Range("D26").Value = "0.0₃0985"
$ws.Range("E26").Value = "  -6.91%  "

$ws.Range("D27").Value = "2.587.31"
$ws.Range("E27").Value = "  -2.25%  "

$ws.Range("E28").Value = "  -3.13%  "

$ws.Range("B29").Value = "Bittensor"
$ws.Range("C29").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D29").Value = "537.95"
$ws.Range("E29").Value = "  -2.19%  "

$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").Value = "0.996"
$ws.Range("E30").Value = "  -0.57%  "

$ws.Range("B31").Value = "Aptos"
$ws.Range("C31").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D31").Value = "7.87"
$ws.Range("E31").Value = "  +0.23%  "

$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "8.29"
$ws.Range("E32").Value = "  -4.06%  "

$ws.Range("D33").Value = "0.147"
$ws.Range("E33").Value = "  -5.04%  "

$ws.Range("D34").Value = "1.88"
$ws.Range("E34").Value = "  -2.75%  "

$ws.Range("D35").Value = "1.64"
$ws.Range("E35").Value = "  +1.48%  "

$ws.Range("D36").Value = "5.74"
$ws.Range("E36").Value = "  -5.32%  "

$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  +0.02%  "

$ws.Range("D38").Value = "4.81"
$ws.Range("E38").Value = "  -2.77%  "

$ws.Range("E39").Value = "  -0.29%  "

$ws.Range("D40").Value = "18.25"
$ws.Range("E40").Value = "  -3.13%  "

$ws.Range("E41").Value = "  +2.47%  "

$ws.Range("D42").Value = "140.02"
$ws.Range("E42").Value = "  -5.57%  "

$ws.Range("E43").Value = "  +0.07%  "

$ws.Range("D44").Value = "40.41"
$ws.Range("E44").Value = "  -1.23%  "

$ws.Range("D45").Value = "2.31"
$ws.Range("E45").Value = "  -3.06%  "

$ws.Range("D46").Value = "143.89"
$ws.Range("E46").Value = "  -4.42%  "

$ws.Range("D47").Value = "3.62"
$ws.Range("E47").Value = "  -1.77%  "

$ws.Range("D48").Value = "21.56"
$ws.Range("E48").Value = "  -1.28%  "

$ws.Range("D49").Value = "0.0530"
$ws.Range("E49").Value = "  -3.39%  "

$ws.Range("D50").Value = "0.592"
$ws.Range("E50").Value = "  -0.96%  "

$ws.Range("E51").Value = "  -2.87%  "

# Restore the default style on the touched cells so only the values differ from before.
$textRange.Style = "Normal"
